# "this added last report 01-03-25"
# Update the daily cash-denomination figures on Sheet1 with the latest
# counted quantities. Dependent formulas (totals, row products, grand
# totals) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Top block (rows 3-9): "Quentity" columns K (J block) and P (O block) ---
$ws.Range("K3").Value = 17
$ws.Range("K4").Value = 22
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 100
$ws.Range("K9").Value = 6

$ws.Range("P3").Value = 19
$ws.Range("P4").Value = 31
$ws.Range("P5").Value = 3
$ws.Range("P6").Value = 26
$ws.Range("P7").Value = 11
$ws.Range("P8").ClearContents()

# --- RSO blocks (rows 18-26): quantity columns H, M, R ---
$ws.Range("H18").Value = 44
$ws.Range("H19").Value = 144
$ws.Range("H20").Value = 39
$ws.Range("H21").Value = 126
$ws.Range("H22").Value = 32
$ws.Range("H24").ClearContents()

$ws.Range("M18").Value = 4
$ws.Range("M19").Value = 36
$ws.Range("M20").Value = 6
$ws.Range("M21").Value = 32
$ws.Range("M22").Value = 40
$ws.Range("M23").Value = 36
$ws.Range("M24").Value = 8
$ws.Range("M26").Value = 1

$ws.Range("R18").Value = 16
$ws.Range("R19").Value = 42
$ws.Range("R20").Value = 17
$ws.Range("R21").Value = 31
$ws.Range("R22").Value = 45
$ws.Range("R23").Value = 18
$ws.Range("R24").Value = 23

# --- Restore the last active selection shown in the file ---
$ws.Range("R25").Select()
